# Sheet1!G6 gets the text "kjk" (stored as a shared string, as Excel does
# for plain text cell entry), and becomes the active/selected cell - matching
# the target worksheet's <sheetData>/<selection> state.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G6").Value = "kjk"
$ws.Range("G6").Select() | Out-Null
